$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Phase 1: create all new paragraphs (style + plain concatenated text)
# Anchor: the last paragraph in the body is "BILAGA 1 - Fridlysta arter"
# (Title style), right before the sectPr. New "Knärot" paragraphs go after it.
# ------------------------------------------------------------------
$anchor = $d.Paragraphs.Last

# --- New paragraph 0: style=Heading1 ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = 'Heading1'
$pStart0 = $newp.Range.Start
$d.Range($pStart0, $pStart0).InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$anchor = $newp

# --- New paragraph 1: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart1 = $newp.Range.Start
$d.Range($pStart1, $pStart1).InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$anchor = $newp

# --- New paragraph 2: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart2 = $newp.Range.Start
$d.Range($pStart2, $pStart2).InsertAfter('Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$anchor = $newp

# --- New paragraph 3: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart3 = $newp.Range.Start
$d.Range($pStart3, $pStart3).InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$anchor = $newp

# --- New paragraph 4: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart4 = $newp.Range.Start
$d.Range($pStart4, $pStart4).InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$anchor = $newp

# --- New paragraph 5: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart5 = $newp.Range.Start
$d.Range($pStart5, $pStart5).InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$anchor = $newp

# --- New paragraph 6: style=Heading2 ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = 'Heading2'
$pStart6 = $newp.Range.Start
$d.Range($pStart6, $pStart6).InsertAfter('Referenser - knärot')
$anchor = $newp

# --- New paragraph 7: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart7 = $newp.Range.Start
$d.Range($pStart7, $pStart7).InsertAfter('de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025')
$anchor = $newp

# --- New paragraph 8: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart8 = $newp.Range.Start
$d.Range($pStart8, $pStart8).InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 ')
$anchor = $newp

# --- New paragraph 9: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart9 = $newp.Range.Start
$d.Range($pStart9, $pStart9).InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853')
$anchor = $newp

# --- New paragraph 10: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart10 = $newp.Range.Start
$d.Range($pStart10, $pStart10).InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.')
$anchor = $newp

# --- New paragraph 11: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart11 = $newp.Range.Start
$d.Range($pStart11, $pStart11).InsertAfter('Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$anchor = $newp

# --- New paragraph 12: style=None ---
$anchor.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$newp.Style = "Normal"
$pStart12 = $newp.Range.Start
$d.Range($pStart12, $pStart12).InsertAfter('SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala ')
$anchor = $newp

# ------------------------------------------------------------------
# Phase 2: apply run-level character formatting (italic/bold) now that
# all paragraph boundaries are finalized. Doing this interleaved with
# paragraph creation causes the COM host to "leak" the last-applied
# character formatting into later InsertParagraphAfter()/InsertAfter()
# calls, so all formatting is deferred to this second pass.
# ------------------------------------------------------------------
$d.Range($pStart2 + 34, $pStart2 + 116).Font.Italic = $true
$d.Range($pStart2 + 278, $pStart2 + 483).Font.Italic = $true
$d.Range($pStart2 + 490, $pStart2 + 608).Font.Italic = $true
$d.Range($pStart3 + 205, $pStart3 + 1070).Font.Italic = $true
$d.Range($pStart7 + 33, $pStart7 + 113).Font.Italic = $true
$d.Range($pStart8 + 62, $pStart8 + 176).Font.Italic = $true
$d.Range($pStart9 + 117, $pStart9 + 207).Font.Italic = $true
$d.Range($pStart10 + 54, $pStart10 + 121).Font.Italic = $true
$d.Range($pStart11 + 22, $pStart11 + 57).Font.Italic = $true
$d.Range($pStart12 + 25, $pStart12 + 61).Font.Italic = $true

# --- Update date in first-page header (wdHeaderFooterFirstPage = 2) ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$null = $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
